$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/legal-document-source"

# Version
$meta.Range("B3").Value = "8.0.0"

# Date
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet "Include from Legal Document S" (CodeSystem info) ---
$codes = $wb.Worksheets.Item("Include from Legal Document S")

# System URI
$codes.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/legal-document-source"
